$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The s1Protocol column (H) holds the kit catalog number "E7760" for every
# data row (2-49). The commit updates that catalog number to "E7420".
$ws.Range("H2:H49").Value = "E7420"
